$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 5 "left'" trial rows (2-6) and 5 "right'" trial
# rows (7-11), each with 4 redundant repeats. Trim each group down to its
# first 2 rows (header + one "None" row + one real-answer row), leaving a
# clean 5-row table (header + 2 left rows + 2 right rows).

# 1) Remove the 3 extra duplicate "left'" rows (old rows 4-6). Rows 7-11
#    shift up to become rows 4-8.
$ws.Rows("4:6").Delete() | Out-Null

# 2) Remove the 3 extra duplicate "right'" rows, which now sit at rows 6-8.
$ws.Rows("6:8").Delete() | Out-Null

# Leave the same rows selected, as was left behind after the deletion.
$ws.Range("A6:XFD8").Select() | Out-Null
